$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz2")

# --- Table 1 (rows 37-42) ---
# Row 37: change allocation split from 60/40 to 80/20
$ws.Range("F37").Value = 0.8
$ws.Range("G37").Value = 0.2

# Row 38: change allocation split, add a 1000 contribution
$ws.Range("F38").Value = 0.8
$ws.Range("G38").Value = 0.2
$ws.Range("H38").Value = 1000

# Row 39: change allocation split
$ws.Range("F39").Value = 0.8
$ws.Range("G39").Value = 0.2

# Row 42: change allocation split, increase withdrawal
$ws.Range("F42").Value = 0.5
$ws.Range("G42").Value = 0.5
$ws.Range("I42").Value = 15000

# --- Table 2 (rows 48-53) ---
# Row 48: change allocation split
$ws.Range("F48").Value = 0.8
$ws.Range("G48").Value = 0.2

# Row 49: change allocation split, add a 1000 contribution
$ws.Range("F49").Value = 0.8
$ws.Range("G49").Value = 0.2
$ws.Range("H49").Value = 1000

# Row 50: change allocation split
$ws.Range("F50").Value = 0.8
$ws.Range("G50").Value = 0.2

# Row 51: change allocation split
$ws.Range("F51").Value = 0.7
$ws.Range("G51").Value = 0.3

# Row 53: change allocation split, clear the withdrawal cell
$ws.Range("F53").Value = 0.5
$ws.Range("G53").Value = 0.5
$ws.Range("I53").NumberFormat = "0.00"
$ws.Range("I53").ClearContents()

# --- Selection / view state ---
$ws.Range("J53").Select()
